# Course Project 2: Plot 1 Timings
# Add the "Macbook Pro" test-run results block (rows 17-21) and bold the
# "HP Omen" section header (A11) to match the other section headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Make the "HP Omen" header bold, like the other section headers ---
$ws.Range("A11").Font.Bold = $true

# --- New section header: Macbook Pro ---
$ws.Range("A17").Value = "Macbook Pro (2.6Ghz i5, 8Gb RAM)"
$ws.Range("A17").Font.Bold = $true

# --- Macbook Pro timing rows ---
$ws.Range("B18").Value = "Read NEI Data"
$ws.Range("C18").Value = 23.32
$ws.Range("D18").Formula = "=C18"

$ws.Range("B19").Value = "Read Classification Codes"
$ws.Range("C19").Value = 23.37
$ws.Range("D19").Formula = "=C19-C18"

$ws.Range("B20").Value = "Aggregate Data"
$ws.Range("C20").Value = 33.65
$ws.Range("D20").Formula = "=C20-C19"

$ws.Range("B21").Value = "Create Bar Plot"
$ws.Range("C21").Value = 33.74
$ws.Range("D21").Formula = "=C21-C20"

# --- View tweaks: zoom in a bit and move the selection ---
$ws.Range("O27").Select() | Out-Null
$excel.ActiveWindow.Zoom = 120
